$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (header years): add 2021, 2022, 2023 with the same style as P3/Q3 (s=6)
$ws.Range("R3").Value = 2021
$ws.Range("S3").Value = 2022
$ws.Range("T3").Value = 2023
$ws.Range("P3:Q3").Copy()
$ws.Range("R3:T3").PasteSpecial(-4122)  # xlPasteFormats

# Row 4 (data row, total appeals): add 4301, 3690, 2620 with same style as P4/Q4 (s=7)
$ws.Range("R4").Value = 4301
$ws.Range("S4").Value = 3690
$ws.Range("T4").Value = 2620
$ws.Range("P4:Q4").Copy()
$ws.Range("R4:T4").PasteSpecial(-4122)

# Row 5 (data row, positively resolved): add 427, 280, 264 with same style as P5/Q5 (s=9)
$ws.Range("R5").Value = 427
$ws.Range("S5").Value = 280
$ws.Range("T5").Value = 264
$ws.Range("P5:Q5").Copy()
$ws.Range("R5:T5").PasteSpecial(-4122)

# Row 2 (thick bottom border spacer row): extend empty bordered cells to R2:T2 (style s=3)
$ws.Range("Q2").Copy()
$ws.Range("R2:T2").PasteSpecial(-4122)

$excel.CutCopyMode = 0
